$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1669288.5
$ws.Range("I19").Value = 9999999
$ws.Range("J19").Value = 3146.4
$ws.Range("K19").Value = 9999999
$ws.Range("L19").Value = 3146.4
$ws.Range("M19").Value = -9999824
$ws.Range("N19").Value = -3496.4
$ws.Range("H116").Value = 18114.428
$ws.Range("I116").Value = 100000
$ws.Range("J116").Value = 4466.8335
$ws.Range("K116").Value = 100000
$ws.Range("L116").Value = 4466.8335
$ws.Range("M116").Value = -96558
$ws.Range("N116").Value = -11350.8335
$ws.Range("H125").Value = 1140
$ws.Range("I125").Value = 1140
$ws.Range("K125").Value = 10260
$ws.Range("M125").Value = -7800
$ws.Range("H135").Value = 234.25
$ws.Range("I135").Value = 114.13333
$ws.Range("J135").Value = 2036
$ws.Range("K135").Value = 1027.19997
$ws.Range("L135").Value = 18324
$ws.Range("M135").Value = 1507.80003
$ws.Range("N135").Value = -23394
$ws.Range("H139").Value = 69971.664
$ws.Range("J139").Value = 69971.664
$ws.Range("L139").Value = 69971.664
$ws.Range("N139").Value = -80251.664
$ws.Range("H140").Value = 75542.31
$ws.Range("J140").Value = 75542.31
$ws.Range("L140").Value = 75542.31
$ws.Range("N140").Value = -85902.31

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6371.091
$ws.Range("I32").Value = 4001.6216
$ws.Range("K32").Value = 4001.6216
$ws.Range("M32").Value = -3714.6216
$ws.Range("H63").Value = 1877.75
$ws.Range("I63").Value = 1870.3334
$ws.Range("K63").Value = 1870.3334
$ws.Range("M63").Value = -1184.3334
$ws.Range("H66").Value = 1877.75
$ws.Range("I66").Value = 1870.3334
$ws.Range("K66").Value = 9351.666999999999
$ws.Range("M66").Value = -5919.666999999999
$ws.Range("H109").Value = 58655.5
$ws.Range("J109").Value = 58655.5
$ws.Range("L109").Value = 58655.5
$ws.Range("N109").Value = -61429.5

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3757.4
$ws.Range("I62").Value = 3443.5
$ws.Range("J62").Value = 3966.6667
$ws.Range("K62").Value = 3443.5
$ws.Range("L62").Value = 3966.6667
$ws.Range("M62").Value = -2819.5
$ws.Range("N62").Value = -5214.6667
$ws.Range("H65").Value = 3757.4
$ws.Range("I65").Value = 3443.5
$ws.Range("J65").Value = 3966.6667
$ws.Range("K65").Value = 17217.5
$ws.Range("L65").Value = 19833.3335
$ws.Range("M65").Value = -14097.5
$ws.Range("N65").Value = -26073.3335
$ws.Range("H70").Value = 42499.75
$ws.Range("J70").Value = 42499.75
$ws.Range("L70").Value = 42499.75
$ws.Range("N70").Value = -43129.75
$ws.Range("H73").Value = 42499.75
$ws.Range("J73").Value = 42499.75
$ws.Range("L73").Value = 42499.75
$ws.Range("N73").Value = -44683.75
$ws.Range("H107").Value = 1078.3572
$ws.Range("I107").Value = 985.5714
$ws.Range("K107").Value = 985.5714
$ws.Range("M107").Value = 934.4286
$ws.Range("H134").Value = 3210.3572
$ws.Range("I134").Value = 2931.9
$ws.Range("K134").Value = 8795.700000000001
$ws.Range("M134").Value = -6260.700000000001

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 392.2857
$ws.Range("I5").Value = 335.55554
$ws.Range("J5").Value = 732.6667
$ws.Range("K5").Value = 1006.66662
$ws.Range("L5").Value = 2198.0001
$ws.Range("M5").Value = -894.66662
$ws.Range("N5").Value = -2422.0001
$ws.Range("H12").Value = 94.07692
$ws.Range("I12").Value = 63.833332
$ws.Range("J12").Value = 120
$ws.Range("K12").Value = 191.499996
$ws.Range("L12").Value = 360
$ws.Range("M12").Value = -18.49999600000001
$ws.Range("N12").Value = -706
$ws.Range("H135").Value = 392.2857
$ws.Range("I135").Value = 335.55554
$ws.Range("J135").Value = 732.6667
$ws.Range("K135").Value = 3019.99986
$ws.Range("L135").Value = 6594.0003
$ws.Range("M135").Value = -484.9998599999999
$ws.Range("N135").Value = -11664.0003

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 22657
$ws.Range("J26").Value = 22657
$ws.Range("L26").Value = 22657
$ws.Range("N26").Value = -23217
$ws.Range("H50").Value = 22657
$ws.Range("J50").Value = 22657
$ws.Range("L50").Value = 22657
$ws.Range("N50").Value = -23653
$ws.Range("H107").Value = 110
$ws.Range("I107").Value = 96.666664
$ws.Range("K107").Value = 96.666664
$ws.Range("M107").Value = 1823.333336
$ws.Range("H113").Value = 1437.125
$ws.Range("I113").Value = 1181.5
$ws.Range("J113").Value = 1522.3334
$ws.Range("K113").Value = 1181.5
$ws.Range("L113").Value = 1522.3334
$ws.Range("M113").Value = 988.5
$ws.Range("N113").Value = -5862.3334
$ws.Range("H132").Value = 2139549
$ws.Range("I132").Value = 3206832.2
$ws.Range("J132").Value = 4982.6665
$ws.Range("K132").Value = 9620496.600000001
$ws.Range("L132").Value = 14947.9995
$ws.Range("M132").Value = -9617966.600000001
$ws.Range("N132").Value = -20007.9995

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1827.8387
$ws.Range("I7").Value = 1729.7587
$ws.Range("J7").Value = 3250
$ws.Range("K7").Value = 1729.7587
$ws.Range("L7").Value = 3250
$ws.Range("M7").Value = -1617.7587
$ws.Range("N7").Value = -3474
$ws.Range("H12").Value = 1002500
$ws.Range("J12").Value = 5000
$ws.Range("L12").Value = 5000
$ws.Range("N12").Value = -5340
$ws.Range("H61").Value = 2646.2666
$ws.Range("I61").Value = 2454.889
$ws.Range("J61").Value = 2933.3333
$ws.Range("K61").Value = 2454.889
$ws.Range("L61").Value = 2933.3333
$ws.Range("M61").Value = -2252.889
$ws.Range("N61").Value = -3337.3333
$ws.Range("H113").Value = 2646.2666
$ws.Range("I113").Value = 2454.889
$ws.Range("J113").Value = 2933.3333
$ws.Range("K113").Value = 2454.889
$ws.Range("L113").Value = 2933.3333
$ws.Range("M113").Value = -284.8890000000001
$ws.Range("N113").Value = -7273.3333
$ws.Range("H122").Value = 9810.1
$ws.Range("I122").Value = 9683.5
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 29050.5
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -26600.5
$ws.Range("N122").Value = -34900
$ws.Range("H126").Value = 1827.8387
$ws.Range("I126").Value = 1729.7587
$ws.Range("J126").Value = 3250
$ws.Range("K126").Value = 5189.2761
$ws.Range("L126").Value = 9750
$ws.Range("M126").Value = -2719.2761
$ws.Range("N126").Value = -14690
$ws.Range("H136").Value = 3574
$ws.Range("I136").Value = 3763.6667
$ws.Range("J136").Value = 3005
$ws.Range("K136").Value = 11291.0001
$ws.Range("L136").Value = 9015
$ws.Range("M136").Value = -8741.000100000001
$ws.Range("N136").Value = -14115

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 11280
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 11280
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 11280
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -12320
$ws.Range("H113").Value = 528.8077
$ws.Range("I113").Value = 349.30768
$ws.Range("J113").Value = 708.3077
$ws.Range("K113").Value = 1047.92304
$ws.Range("L113").Value = 2124.9231
$ws.Range("M113").Value = 1122.07696
$ws.Range("N113").Value = -6464.9231
$ws.Range("H132").Value = 1182.4219
$ws.Range("I132").Value = 875.45654
$ws.Range("J132").Value = 1966.8889
$ws.Range("K132").Value = 2626.36962
$ws.Range("L132").Value = 5900.6667
$ws.Range("M132").Value = -96.36961999999994
$ws.Range("N132").Value = -10960.6667
